$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 21.95262029704519

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 6.189590430959694
